# trendData: add the 2024-11-04 rows for participants 111/112/113 (9 new
# rows of vitals-trend data), extending the existing A1:G10 table down to
# A1:G19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append, in sheet order (row number + the 7 column values:
# Date, Participant Number, Block Name, Probe number, Vital Sign, Trend,
# Timestamp).
$newRows = @(
    @{ Row = 11; Date = "2024-11-04"; Participant = "111"; Block = "AA111";  Probe = 4; Vital = "Heart Rate";        Trend = "Increasing"; Timestamp = "2024-11-04T06:54:28.223" },
    @{ Row = 12; Date = "2024-11-04"; Participant = "111"; Block = "AA111";  Probe = 4; Vital = "Blood Pressure";    Trend = "Decreasing"; Timestamp = "2024-11-04T06:54:44.224" },
    @{ Row = 13; Date = "2024-11-04"; Participant = "111"; Block = "AA111";  Probe = 4; Vital = "Oxygen Saturation"; Trend = "Decreasing"; Timestamp = "2024-11-04T06:55:08.077" },
    @{ Row = 14; Date = "2024-11-04"; Participant = "112"; Block = "AA112";  Probe = 5; Vital = "Blood Pressure";    Trend = "Decreasing"; Timestamp = "2024-11-04T06:58:10.855" },
    @{ Row = 15; Date = "2024-11-04"; Participant = "112"; Block = "AA112";  Probe = 5; Vital = "Heart Rate";        Trend = "Increasing"; Timestamp = "2024-11-04T06:58:17.154" },
    @{ Row = 16; Date = "2024-11-04"; Participant = "112"; Block = "AA112";  Probe = 5; Vital = "Oxygen Saturation"; Trend = "Increasing"; Timestamp = "2024-11-04T06:58:17.772" },
    @{ Row = 17; Date = "2024-11-04"; Participant = "113"; Block = "Aa1123"; Probe = 6; Vital = "Heart Rate";        Trend = "Increasing"; Timestamp = "2024-11-04T07:02:43.104" },
    @{ Row = 18; Date = "2024-11-04"; Participant = "113"; Block = "Aa1123"; Probe = 6; Vital = "Oxygen Saturation"; Trend = "Static";      Timestamp = "2024-11-04T07:02:47.245" },
    @{ Row = 19; Date = "2024-11-04"; Participant = "113"; Block = "Aa1123"; Probe = 6; Vital = "Blood Pressure";    Trend = "Decreasing"; Timestamp = "2024-11-04T07:02:55.084" }
)

foreach ($r in $newRows) {
    $i = $r.Row

    # Columns A (Date) and B (Participant Number) hold values that look like
    # a date / a plain number ("2024-11-04", "111", ...), but the source data
    # keeps them as literal text (same as the existing A2:B10 cells). Mark the
    # cells as Text *before* assigning so Excel doesn't auto-convert them to a
    # date serial / number.
    $ws.Range("A$i").NumberFormat = "@"
    $ws.Range("A$i").Value = $r.Date
    $ws.Range("B$i").NumberFormat = "@"
    $ws.Range("B$i").Value = $r.Participant

    # Column C (Block Name) is alphanumeric text - no conversion risk.
    $ws.Range("C$i").Value = $r.Block

    # Column D (Probe number) is a genuine number, same as the existing rows.
    $ws.Range("D$i").Value = $r.Probe

    # Columns E/F (Vital Sign / Trend) are plain words.
    $ws.Range("E$i").Value = $r.Vital
    $ws.Range("F$i").Value = $r.Trend

    # Column G (Timestamp) contains a literal "T" so Excel won't parse it as
    # a date/time value - stays text on its own.
    $ws.Range("G$i").Value = $r.Timestamp
}

# The workbook ignores the "number stored as text" warning over the table;
# extend that suppression to the newly-added rows as well.
try {
    $ws.Range("A1:G19").Errors.Item(9).Ignore = $true
} catch {
    # Older/partial COM surfaces may not implement Range.Errors - harmless
    # to skip, the data itself is unaffected.
}
